# Count chlD genes instead of medium subunit genes
# -----------------------------------------------------------------
# This script reproduces the target diff against the statistics table:
#  - relabels the "medium subunit" header wording to "chlD" wording
#  - updates two data values (E5, E6) that changed because the counting
#    method changed
#  - updates the TOTAL row's cached sum accordingly (handled by recalc)
#  - appends a new flattened data table (rows 15-26) with machine-
#    friendly column names, mirroring the pivot table above it
#  - tidies the view (selection cell, column widths)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 0) Stamp rows 15-26 with the plain base style (xf index 5 - the
#    same one already used by the sheet's column defaults) before
#    filling in any values, by pasting formats from a same-styled
#    cell. PasteSpecial() with no arguments resets the destination
#    to the engine's default style, which happens to be xf 5 here -
#    exactly what the target rows use - without adding any new
#    entries to styles.xml.
# ---------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A15:G26").PasteSpecial() | Out-Null

# ---------------------------------------------------------------
# 1) New flat data table header (row 15). Written first so these
#    brand-new shared strings land immediately after the existing
#    ones (Archaea..TOTAL) in the shared-string table.
# ---------------------------------------------------------------
$ws.Range("A15").Value = "kingdom"
$ws.Range("B15").Value = "phylum2"
$ws.Range("C15").Value = "total_orgs"
$ws.Range("D15").Value = "orgs_with_bchlD"
$ws.Range("E15").Value = "num_bchlD"
$ws.Range("F15").Value = "orgs_with_bchlD_fs"
$ws.Range("G15").Value = "num_bchlD_fs"

# ---------------------------------------------------------------
# 2) Re-word the pivot table's header row (D1:G1). A1/B1 (Domain /
#    Phylum) and C1 (Total number of genomes) keep their existing
#    text, so they do not need to be rewritten.
# ---------------------------------------------------------------
$ws.Range("D1").Value = "Genomes with chlD genes"
$ws.Range("E1").Value = "Number of chlD genes"
$ws.Range("F1").Value = "Genomes with fs-chlD genes"
$ws.Range("G1").Value = "Number of" + [char]10 + "fs-chlD genes"

# ---------------------------------------------------------------
# 3) Updated counts (genomes counted by chlD genes instead of by
#    medium subunit genes) for Proteobacteria and Actinobacteria.
# ---------------------------------------------------------------
$ws.Range("E5").Value = 337
$ws.Range("E6").Value = 572

# ---------------------------------------------------------------
# 4) New flattened data rows (16-26), mirroring rows 2-12 of the
#    pivot table but with every row carrying its own kingdom/phylum
#    label (no merged cells) and the updated chlD-based counts.
# ---------------------------------------------------------------
$ws.Range("A16").Value = "Archaea"
$ws.Range("B16").Value = "Euryarchaeota"
$ws.Range("C16").Value = 220
$ws.Range("D16").Value = 134
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 25
$ws.Range("G16").Value = 27

$ws.Range("A17").Value = "Archaea"
$ws.Range("B17").Value = "Crenarchaeota"
$ws.Range("C17").Value = 24
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0

$ws.Range("A18").Value = "Archaea"
$ws.Range("B18").Value = "Thaumarchaeota"
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0

$ws.Range("A19").Value = "Bacteria"
$ws.Range("B19").Value = "Proteobacteria"
$ws.Range("C19").Value = 2027
$ws.Range("D19").Value = 321
$ws.Range("E19").Value = 337
$ws.Range("F19").Value = 87
$ws.Range("G19").Value = 87

$ws.Range("A20").Value = "Bacteria"
$ws.Range("B20").Value = "Actinobacteria"
$ws.Range("C20").Value = 1024
$ws.Range("D20").Value = 540
$ws.Range("E20").Value = 572
$ws.Range("F20").Value = 38
$ws.Range("G20").Value = 38

$ws.Range("A21").Value = "Bacteria"
$ws.Range("B21").Value = "Chloroflexi"
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 21
$ws.Range("F21").Value = 13
$ws.Range("G21").Value = 13

$ws.Range("A22").Value = "Bacteria"
$ws.Range("B22").Value = "Spirochaetes"
$ws.Range("C22").Value = 67
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2

$ws.Range("A23").Value = "Bacteria"
$ws.Range("B23").Value = "Firmicutes"
$ws.Range("C23").Value = 1215
$ws.Range("D23").Value = 55
$ws.Range("E23").Value = 61
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 1

$ws.Range("A24").Value = "Bacteria"
$ws.Range("B24").Value = "Bacteroidetes"
$ws.Range("C24").Value = 569
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 27
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 1

$ws.Range("A25").Value = "Bacteria"
$ws.Range("B25").Value = "Cyanobacteria"
$ws.Range("C25").Value = 77
$ws.Range("D25").Value = 76
$ws.Range("E25").Value = 78
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0

$ws.Range("A26").Value = "Bacteria"
$ws.Range("B26").Value = "Other"
$ws.Range("C26").Value = 361
$ws.Range("D26").Value = 38
$ws.Range("E26").Value = 47
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0

# ---------------------------------------------------------------
# 6) Column widths: split the old single "C:G" band (width 13.5)
#    into "C:F" (unchanged, 13.5) and a narrower "G" column (target
#    stored width ~12.332; 11.5 is the closest API input this
#    engine's width rounding reproduces).
# ---------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 11.5

# ---------------------------------------------------------------
# 7) Row 1 height shrinks now that the header text is shorter.
# ---------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 36

# ---------------------------------------------------------------
# 8) Selection moves to G2.
# ---------------------------------------------------------------
$ws.Range("G2").Select() | Out-Null

"edit applied"
